$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C2:C6) from serial 45183 to 45184 (one day later)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45184
}
